# Apply updated symbol list values (coin name/link/price/volume changes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'7.28%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'8.05%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.303"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'3.62%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'3.81%"
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = "'KuCoinToken"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'8.653"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.34%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "'FTXToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'1.927"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.65%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "'BTSEToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'2.977"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.42%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'MXToken"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.9421"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'2.07%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.1361"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'23.06%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'WazirX"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.1970"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.04%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.09198"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.80%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.03553"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'5.87%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.09594"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.08%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.001330"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-3.98%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.006064"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.99%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'HotbitToken"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.004320"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.70%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.366"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-2.07%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'GateToken"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'4.522"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.83%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'BitpandaEcosystemToken"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.3511"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.62%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'MCDex"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'7.226"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'15.52%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'ProBitToken"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.1333"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.59%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'ZBToken"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.2562"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'5.19%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'CoinExToken"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'0.04428"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.95%"
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'BitKan"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'0.001223"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.29%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.80%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.07%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02486"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'15.75%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05201"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.90%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007649"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.59%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'5.80%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.009132"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'5.16%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002162"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'8.53%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009915"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'13.16%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006659"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.46%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.11%"
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'CoinbaseStockToken"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.002403"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'139.76%"
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'BOLO"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.003345"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-1.16%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.11%"
$ws.Range("E51").Style = "Normal"
